$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B9").Value = 44172
$ws.Range("B9").NumberFormat = "d-mmm"

$ws.Range("C9").Value = "Meeting + Webseite bauen (Slider + weiterer Div Bereich) "
$ws.Range("D9").Value = "14:00 Uhr "
$ws.Range("E9").Value = "17:00 Uhr"
$ws.Range("F9").Value = 3

$ws.Range("G9").Select()
